$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that change per row: D (4), J (10), K (11), L (12), M (13), P (16)
$cols = @(4, 10, 11, 12, 13, 16)

# Capture the original ("before") values for each data row (2..19)
$orig = @{}
for ($r = 2; $r -le 19; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $orig[$r] = $rowVals
}

# Mapping: new row -> old row whose original values it should receive
$mapping = @{
    2  = 6
    3  = 12
    4  = 17
    5  = 18
    6  = 7
    7  = 15
    8  = 11
    9  = 16
    10 = 8
    11 = 4
    12 = 5
    13 = 9
    14 = 13
    15 = 19
    16 = 2
    17 = 14
    18 = 10
    19 = 3
}

foreach ($newRow in $mapping.Keys) {
    $srcRow = $mapping[$newRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($newRow, $c).Value = $orig[$srcRow][$c]
    }
}
